$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Recommandations")
$data1 = @(
  @("Titre", "Jours en Hausse", "Jours en Baisse", "Variation Totale (%)", "Dernière Variation (%)", "Recommandation", "Stratégie"),
  @("BRVM - SERVICES PUBLICS", 0, 8, 3429.75, 112.71, "🟡 Observer", "➖ Neutre"),
  @("AIR LIQUIDE CI", 0, 4, 2800, 675, "🟡 Observer", "➖ Neutre"),
  @("NEI-CEDA CI", 0, 4, 2785, 690, "🟡 Observer", "➖ Neutre"),
  @("BRVM - AUTRES SECTEURS", 0, 4, 2447.94, 618.5599999999999, "🟡 Observer", "➖ Neutre"),
  @("BRVM - DISTRIBUTION", 0, 4, 2083.56, 528.0700000000001, "🟡 Observer", "➖ Neutre"),
  @("BRVM - AGRICULTURE", 0, 4, 1463.82, 368.63, "🟡 Observer", "➖ Neutre"),
  @("BRVM - TRANSPORT", 0, 4, 1447.83, 354.92, "🟡 Observer", "➖ Neutre"),
  @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 753.2, 192.79, "🟡 Observer", "➖ Neutre"),
  @("BRVM - CONSOMMATION DE BASE         (**)", 0, 3, 570.23, 190.31, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRESTIGE", 0, 4, 564.17, 141.16, "🟡 Observer", "➖ Neutre"),
  @("BRVM - FINANCES", 0, 4, 564.1, 141.63, "🟡 Observer", "➖ Neutre"),
  @("BRVM - SERVICES FINANCIERS", 0, 4, 554.4, 139.2, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIELS", 0, 4, 513, 126.65, "🟡 Observer", "➖ Neutre"),
  @("BRVM - ENERGIE", 0, 4, 433.64, 108.73, "🟡 Observer", "➖ Neutre"),
  @("BRVM - TELECOMMUNICATIONS", 0, 4, 388.75, 96.72, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIE              (**)", 0, 1, 219.1, 219.1, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIE                 (**)", 0, 1, 215.1, 215.1, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIE         (**)", 0, 1, 214.39, 214.39, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIE                (**)", 0, 1, 214.18, 214.18, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRINCIPAL            (**)", 0, 1, 205.89, 205.89, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRINCIPAL                   (**)", 0, 1, 205.64, 205.64, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRINCIPAL                 (**)", 0, 1, 204.98, 204.98, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRINCIPAL                    (**)", 0, 1, 204.3, 204.3, "🟡 Observer", "➖ Neutre"),
  @("BRVM - CONSOMMATION DE BASE           (**)", 0, 1, 193.79, 193.79, "🟡 Observer", "➖ Neutre"),
  @("CFAO MOTORS CI (CFAC)", 3, 0, 19.72, 5, "🟢 Achat", "✅ Renforcer"),
  @("SAFCA CI (SAFC)", 3, 0, 17.4, 5.96, "🟢 Achat", "✅ Renforcer"),
  @("UNIWAX CI (UNXC)", 2, 0, 10.88, 5, "🟡 Observer", "➖ Neutre"),
  @("CIE CI (CIEC)", 2, 0, 7.43, 2.89, "🟡 Observer", "➖ Neutre"),
  @("SETAO CI (STAC)", 1, 0, 7.27, 7.27, "🟡 Observer", "➖ Neutre"),
  @("ORAGROUP TOGO (ORGT)", 1, 0, 4.55, 4.55, "🟡 Observer", "➖ Neutre"),
  @("TOTALENERGIES MARKETING CI (TTLC)", 1, 0, 4.35, 4.35, "🟡 Observer", "➖ Neutre"),
  @("TRACTAFRIC MOTORS CI (PRSC)", 1, 0, 4.08, 4.08, "🟡 Observer", "➖ Neutre"),
  @("BICI CI (BICC)", 1, 0, 3.97, 3.97, "🟡 Observer", "➖ Neutre"),
  @("ORANGE COTE D'IVOIRE (ORAC)", 1, 1, 0.34, -3.13, "🟡 Observer", "👀 À surveiller"),
  @("TOTAL", 0, 3, 0, 0, "🟡 Observer", "➖ Neutre"),
  @("AIR LIQUIDE CI (SIVC)", 1, 1, -1.88, 4.41, "🟡 Observer", "👀 À surveiller"),
  @("VIVO ENERGY CI (SHEC)", 0, 1, -1.99, -1.99, "🟡 Observer", "➖ Neutre"),
  @("PALM CI (PALC)", 0, 1, -1.99, -1.99, "🟡 Observer", "➖ Neutre"),
  @("SITAB CI (STBC)", 0, 1, -2.02, -2.02, "🟡 Observer", "➖ Neutre"),
  @("BANK OF AFRICA SENEGAL (BOAS)", 0, 1, -2.14, -2.14, "🟡 Observer", "➖ Neutre"),
  @("BANK OF AFRICA NG (BOAN)", 0, 1, -2.5, -2.5, "🟡 Observer", "➖ Neutre"),
  @("AFRICA GLOBAL LOGISTICS CI (SDSC)", 0, 1, -2.68, -2.68, "🟡 Observer", "➖ Neutre"),
  @("SICOR CI (SICC)", 1, 2, -2.7, 2.94, "🟡 Observer", "👀 À surveiller"),
  @("ONATEL BF (ONTBF)", 0, 1, -3.61, -3.61, "🟡 Observer", "➖ Neutre"),
  @("LOTERIE NATIONALE DU BENIN (LNBB)", 0, 1, -3.75, -3.75, "🟡 Observer", "➖ Neutre"),
  @("BERNABE CI (BNBC)", 1, 2, -4.17, 5.51, "🟡 Observer", "👀 À surveiller"),
  @("SERVAIR ABIDJAN CI (ABJC)", 1, 2, -4.67, 7.32, "🟡 Observer", "👀 À surveiller"),
  @("SOLIBRA CI (SLBC)", 0, 1, -7.46, -7.46, "🟡 Observer", "➖ Neutre"),
  @("SICABLE CI (CABC)", 0, 3, -14.93, -4.19, "🔴 Vente", "⚠️ Risque de décrochage")
)

for ($r = 0; $r -lt $data1.Count; $r++) {
  for ($c = 0; $c -lt 7; $c++) {
    $ws1.Cells.Item($r + 1, $c + 1).Value = $data1[$r][$c]
  }
}

$ws2 = $wb.Worksheets.Item("Top_YTD")
$data2 = @(
  @("Titre", "Progression YTD (%)"),
  @("BRVM - SERVICES PUBLICS", 10234808.03),
  @("AIR LIQUIDE CI", 409193),
  @("NEI-CEDA CI", 401662.4),
  @("BRVM - AUTRES SECTEURS", 256806.82),
  @("BRVM - DISTRIBUTION", 148411.78),
  @("BRVM - AGRICULTURE", 47036.27),
  @("BRVM - TRANSPORT", 45434.38),
  @("BRVM - CONSOMMATION DISCRETIONNAIRE", 6804.04),
  @("BRVM-PRESTIGE", 3275.77),
  @("BRVM - FINANCES", 3274.77)
)

for ($r = 0; $r -lt $data2.Count; $r++) {
  for ($c = 0; $c -lt 2; $c++) {
    $ws2.Cells.Item($r + 1, $c + 1).Value = $data2[$r][$c]
  }
}